# Added png for mana and fireball
#
# Insert a new "Cost" column before the CardLayout column (C), which
# shifts CardLayout/Image/Amount one column to the right (C->D, D->E, E->F)
# and pushes the old Cost column (F) to G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Insert()

# New "Cost" header in C1
$ws.Range("C1").Value = "Cost"

# Bring the old Cost values (now shifted into column G) into the new
# column C, keeping them next to Name/Text as in the updated layout.
$ws.Range("C4").Value = $ws.Range("G4").Value2
$ws.Range("C5").Value = $ws.Range("G5").Value2

# Drop the now-duplicated values from the shifted-away old Cost column,
# but keep the (now blank) styled cells in place.
$ws.Range("G1").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("G5").ClearContents()

# The extra blank stub cells that the column insert created in column C
# for rows that never had a Cost value shouldn't stick around.
$ws.Range("C2").Clear()
$ws.Range("C3").Clear()
$ws.Range("C6").Clear()
$ws.Range("C7").Clear()

# New "Amount" data: every card has an amount of 1.
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1

# New card-art image paths (Fireball first, then Mana, to match authoring order).
$ws.Range("E4").Value = "C:\Users\rpswa\WorkSpace\spell-slingers-prototype\assets\fireball.png"
$ws.Range("E2").Value = "C:\Users\rpswa\WorkSpace\spell-slingers-prototype\assets\mana-potion.png"

# Widen the Name/Text/Cost columns slightly and give the new Image column
# enough room to show the long asset paths.
$ws.Columns("B:C").ColumnWidth = 16.5
$ws.Columns("E").ColumnWidth = 61.333333333333336
